$d = $word.ActiveDocument

# Scene 47 ("Lilith"): Pro leaves for school, and Lilith's closing cue on
# that beat becomes an "(exit)" stage direction instead of the reused
# "(neutral neutral)" expression tag. The same cue text also appears
# earlier in the scene attached to a line of spoken dialogue
# ("Lilith (neutral neutral): Oh, it's you.") - that one must stay as-is,
# so match paragraph-by-paragraph on the exact stage-direction-only text
# rather than doing a blind document-wide replace.
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $text = $r.Text.TrimEnd("`r", "`n", [char]7)
    if ($text -eq "Lilith (neutral neutral):") {
        $r.Find.Execute("Lilith (neutral neutral):", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Lilith (exit):", 2)
    }
}
